# Update the Active Candidates list: refresh rows 6-13 with new data
# and remove the now-obsolete rows 14-15 (table shrinks from 15 to 13 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 6 through 13 (columns A-E)
$data = @(
    @(811, "Navan", "ENT AE (Boston + NYC)", "Tom Andrews", "1st Interview"),
    @(829, "Rox", "Sales Engineer (NY / Austin / SF)", "Catherine LaChapelle", "2nd Interview"),
    @(829, "Rox", "Sales Engineer (NY / Austin / SF)", "Matthew Vielkind", "CV Sent"),
    @(829, "Rox", "Sales Engineer (NY / Austin / SF)", "Ryan Kerbs", "3rd Interview"),
    @(830, "Loop", "Customer Operations Manager", "Geoffre Lavy", "3rd Interview"),
    @(830, "Loop", "Customer Operations Manager", "Gilbert Pasquale", "2nd Interview"),
    @(830, "Loop", "Customer Operations Manager", "Pat Ingersoll", "3rd Interview"),
    @(873, "CodeRabbit", "VP Enterprise Sales", "Andrew Sweet", "1st Interview")
)

$row = 6
foreach ($rowData in $data) {
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $row++
}

# Remove the two rows that are no longer needed (old rows 14 and 15)
$ws.Range("A14:E15").Delete()
